$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.501.30"
$ws.Range("E2").Value = "  +0.59%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.913.31"
$ws.Range("E3").Value = "  -0.06%  "

# Row 4
$ws.Range("E4").Value = "  +0.60%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.39"
$ws.Range("E5").Value = "  +0.88%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.006"
$ws.Range("E6").Value = "  +0.43%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4821"
$ws.Range("E7").Value = "  +1.86%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4066"
$ws.Range("E8").Value = "  -0.22%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08144"
$ws.Range("E9").Value = "  +1.35%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.010"
$ws.Range("E10").Value = "  +0.60%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.40"
$ws.Range("E11").Value = "  +4.12%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.899.73"
$ws.Range("E12").Value = "  -0.15%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.995"
$ws.Range("E13").Value = "  +1.73%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.126"
$ws.Range("E14").Value = "  -0.19%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "90.10"
$ws.Range("E15").Value = "  +0.49%  "

# Row 16
$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.008"
$ws.Range("E16").Value = "  +0.56%  "

# Row 17
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.06764"
$ws.Range("E17").Value = "  +1.84%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001040"
$ws.Range("E18").Value = "  +0.93%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.70"
$ws.Range("E19").Value = "  -0.04%  "

# Row 20
$ws.Range("E20").Value = "  +0.51%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "29.516.53"
$ws.Range("E21").Value = "  +0.60%  "

# Row 22
$ws.Range("E22").Value = "  +1.90%  "

# Row 23
$ws.Range("E23").Value = "  +2.67%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.180"
$ws.Range("E24").Value = "  -0.99%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.149.06"
$ws.Range("E25").Value = "  +0.81%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.77"
$ws.Range("E26").Value = "  +0.59%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.380"
$ws.Range("E27").Value = "  +4.99%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.05"
$ws.Range("E28").Value = "  +1.16%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.103"
$ws.Range("E29").Value = "  -0.49%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.79"
$ws.Range("E30").Value = "  +1.73%  "

# Row 31
$ws.Range("E31").Value = "  -4.23%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09527"
$ws.Range("E32").Value = "  -0.13%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.509"
$ws.Range("E33").Value = "  +1.96%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.560"
$ws.Range("E34").Value = "  +0.21%  "

# Row 35
$ws.Range("E35").Value = "  -2.88%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02268"
$ws.Range("E36").Value = "  +0.68%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06089"
$ws.Range("E37").Value = "  +0.00%  "

# Row 38
$ws.Range("E38").Value = "  +0.22%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5951"
$ws.Range("E39").Value = "  +1.02%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.958"
$ws.Range("E40").Value = "  -3.76%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.66"
$ws.Range("E41").Value = "  +5.37%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1854"
$ws.Range("E42").Value = "  +0.81%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.286"
$ws.Range("E43").Value = "  +2.26%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.407"
$ws.Range("E44").Value = "  -5.44%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.59"
$ws.Range("E45").Value = "  +4.18%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07635"
$ws.Range("E46").Value = "  -3.53%  "

# Row 47
$ws.Range("E47").Value = "  +0.52%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.938"
$ws.Range("E48").Value = "  +0.48%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "115.75"
$ws.Range("E49").Value = "  +2.45%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "72.36"
$ws.Range("E50").Value = "  +1.46%  "

# Row 51
$ws.Range("B51").Value = "MXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.410"
$ws.Range("E51").Value = "  +2.47%  "
